$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '87.815.94'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -2.16%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.099.78'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -2.54%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.41'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.45%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '631.06'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +2.56%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.377'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -3.11%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.809'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +17.51%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.999'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.03%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '3.096.21'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -2.50%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.584'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.40%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +1.13%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000243'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -4.26%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.33'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.73%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '87.555.20'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.659.01'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -2.91%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '31.80'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -3.04%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.101.11'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -2.42%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.37'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +3.21%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0000213'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +8.69%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.20'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -1.21%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '421.52'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -3.11%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.36'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -2.37%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.85'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -4.01%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.44'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +6.39%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '83.48'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +11.04%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.32'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -3.03%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.259.07'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -2.71%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.05%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.23%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -8.57%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.11'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -3.46%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.76'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -6.76%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.148'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +15.60%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '501.27'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -5.98%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.77'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -3.45%  '
$ws.Range('B37').NumberFormat = '@'
$ws.Range('B37').Value = 'Fetch.AI'
$ws.Range('C37').NumberFormat = '@'
$ws.Range('C37').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.26'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.74%  '
$ws.Range('B38').NumberFormat = '@'
$ws.Range('B38').Value = 'PancakeSwap'
$ws.Range('C38').NumberFormat = '@'
$ws.Range('C38').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.81'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -2.37%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '22.34'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +1.98%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '22.17'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.61%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.27%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.04%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.364'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -1.68%  '
$ws.Range('B44').NumberFormat = '@'
$ws.Range('B44').Value = 'Stacks'
$ws.Range('C44').NumberFormat = '@'
$ws.Range('C44').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.83'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -4.38%  '
$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value = 'Stellar'
$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.136'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +10.11%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '146.49'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -1.29%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '43.77'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.41%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0655'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +11.40%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '160.86'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -6.53%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.714'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +1.87%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.18'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -3.66%  '
